# "correção nos dados e inicio da analise PNAD 2009"
#
# The sheet had a spurious section-header row ("grandes regiões e unidades
# da federação") inserted at row 6 with no data of its own, which pushed
# every subsequent region's numeric figures one row out of alignment with
# its label (e.g. the "norte" label sat on an empty row while its numbers
# ended up one row down, etc., all the way through "distrito federal").
#
# The fix is simply to delete that stray row 6: every row below it (labels
# together with their data) shifts up by one, realigning each region's
# label with its own figures, and the sheet shrinks from 38 to 37 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Delete()
